# Update Inscritos/Pagos/Homologadas counts in the "Inscricoes" sheet
# as described in the commit diff (values increased/decreased for several rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$updates = @{
    "E7" = 29
    "E10" = 437
    "F10" = 215
    "H10" = 215
    "E11" = 295
    "F11" = 163
    "H11" = 163
    "E12" = 430
    "F12" = 229
    "H12" = 229
    "E13" = 114
    "F13" = 57
    "H13" = 57
    "E15" = 141
    "F15" = 56
    "H15" = 56
    "E16" = 178
    "E17" = 85
    "F17" = 42
    "H17" = 42
    "E22" = 150
    "E23" = 180
    "E25" = 233
    "E26" = 130
    "E27" = 298
    "E28" = 177
    "F28" = 62
    "H28" = 62
    "E29" = 152
    "E34" = 194
    "F34" = 116
    "H34" = 116
    "E35" = 128
    "F35" = 80
    "H35" = 80
    "E36" = 61
    "F36" = 38
    "H36" = 38
    "E37" = 139
    "E41" = 352
    "E42" = 322
    "E43" = 105
    "F43" = 56
    "H43" = 56
    "E44" = 278
    "F44" = 132
    "H44" = 132
    "E47" = 390
    "E48" = 185
    "F48" = 74
    "H48" = 74
    "E49" = 262
    "F49" = 109
    "H49" = 109
    "E50" = 227
    "F50" = 97
    "H50" = 97
    "E51" = 212
    "E52" = 23
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
